$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" summary text on sheet "Hoja1" ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.36 = 12835.61 pesos`n✅ 12835.61 pesos = 3.34 = 935.82 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Update the rate figures on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 297.999
$wsTasas.Range("O10").Value = 3825

$wsTasas.Range("N12").Value = 3844
$wsTasas.Range("O12").Value = 280.26
